$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @("PractitionerRole", "US Core CareTeam Profile", "US Core Organization Profile", "US Core Patient Profile", "US Core Practitioner Profile", "US Core RelatedPerson Profile")

$rows = @(36, 37, 39)

foreach ($r in $rows) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        # Column C is index 3
        $col = 3 + $i
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}
